$d = $word.ActiveDocument

# 1. Remove Dan's "Project scope" backlog notes (the block of paragraphs that
#    runs from "This area is basically ... Business opportunities ..." through
#    "Also, the office administrator will have sole access ..."), which has
#    been transferred out of the Backlog into the SRS document.
$startPar = $null
$endPar = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "This area is basically*") {
        $startPar = $p
    }
    if ($t -like "Also, the office administrator will have sole access*") {
        $endPar = $p
        break
    }
}
if (($startPar -ne $null) -and ($endPar -ne $null)) {
    $delRange = $d.Range($startPar.Range.Start, $endPar.Range.End)
    $delRange.Delete()
}

# 2. The "_GoBack" bookmark used to sit at the end of the "Assumptions and
#    dependencies" paragraph; it now belongs right before the "References"
#    run, at the start of the (now earlier) "References (TBA)" paragraph.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "References*") {
        $target = $d.Range($p.Range.Start, $p.Range.Start)
        $d.Bookmarks.Add("_GoBack", $target)
        break
    }
}

# 3. Drop the stale lastRenderedPageBreak marker in front of "External
#    interface requirements" by re-running Find/Replace over it.
$d.Content.Find.Execute("External interface requirements", $false, $false, $false, $false, $false, $true, 1, $false, "External interface requirements", 2)
